# Apply updates described by the commit: "ji/update debbuging errors codes"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update upper-bound values (column G) for several feature rows
$ws.Range("G6").Value = 400
$ws.Range("G12").Value = 500
$ws.Range("G13").Value = 450
$ws.Range("G14").Value = 1500

# Update the active selection on the sheet to G15
$ws.Range("G15").Select()

# Update the workbook window vertical position (yWindow)
$excel.ActiveWindow.Top = 9105
